$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 47710
$ws.Range("C2").Value = 26717
$ws.Range("D2").Value = 6680
$ws.Range("E2").Value = 14313

$ws.Range("B3").Value = 47710
$ws.Range("C3").Value = 26717
$ws.Range("D3").Value = 6680
$ws.Range("E3").Value = 14313
